$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the text-like values (prices with "." thousand separators, percentages
# with leading/trailing spaces) are written verbatim as text instead of being
# reinterpreted by Excel as numbers/dates, matching the inlineStr cells already
# used throughout this sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.344.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.973.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.37"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -10.43%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.30"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.09%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.849"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.24%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.263.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.76"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -7.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.970.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.187.13"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.70"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.92"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.131"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0646"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.56%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.20"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0961"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.20"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.29%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0212"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.57"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.356.62"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.35"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.07%  "
